# Edit 1: Paragraph 1 heading - remove the hidden _GoBack bookmark while
# keeping the visible text/formatting unchanged (requires rewriting the
# paragraph's runs because the bookmark sits between two of them).
$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End)
$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r w:rsidRPr="00091943"><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>C</w:t></w:r>
            <w:r w:rsidR="00046A17"><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>MP73010</w:t></w:r>
            <w:r w:rsidRPr="00091943"><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> – Assignment 1 exercise</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$r1.InsertXML($xml1)

# Edit 2: Collapse the ">>  >  your" split runs (with the gramStart/gramEnd
# proofing-error markers) into a single plain run reading
# ">>>  your stuff after this line >>>".
$p2 = $d.Paragraphs(4)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>&gt;&gt;&gt;  your stuff after this line &gt;&gt;&gt;</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$r2.InsertXML($xml2)

# Edit 3: "Ben changing things up!" -> "I have successfully modified a
# project on GitHub."
$p3 = $d.Paragraphs(5)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End)
$r3.Text = "I have successfully modified a project on GitHub."

# Edit 4: Fill the first of the two trailing empty paragraphs with three
# separate runs: "Luke" / " Johnstone" / " is a Git!".
$p4 = $d.Paragraphs(6)
$r4 = $d.Range($p4.Range.Start, $p4.Range.Start)
$xml4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>Luke</w:t></w:r>
            <w:r><w:t xml:space="preserve"> Johnstone</w:t></w:r>
            <w:r><w:t xml:space="preserve"> is a Git!</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$r4.InsertXML($xml4)
